$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XAlphaDeals")

# 1) Duplicate row 2 into a new row 3 (same formatting/values), preserving header/row styles.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).PasteSpecial(-4104)  # xlPasteAll

# 2) Update the new row's TestCaseID and Direction to represent the new "sell" deal.
$ws.Range("A3").Value2 = "QA_TestCase_Auto_X-Alpha_005"
$ws.Range("C3").Value2 = "sell"

# 3) Re-apply the boolean cell's format on the new row (keeps same style as N2).
$ws.Range("N2").Copy()
$ws.Range("N3").PasteSpecial(-4122)  # xlPasteFormats

# 4) Change the existing deal's Direction to lowercase "buy".
$ws.Range("C2").Value2 = "buy"

# 5) Autofit the columns that now contain the wider, newly duplicated data.
$ws.Range("D1:D3").EntireColumn.AutoFit()
$ws.Range("E1:E3").EntireColumn.AutoFit()
$ws.Range("F1:F3").EntireColumn.AutoFit()
$ws.Range("H1:H3").EntireColumn.AutoFit()
$ws.Range("I1:I3").EntireColumn.AutoFit()
$ws.Range("J1:J3").EntireColumn.AutoFit()
$ws.Range("K1:K3").EntireColumn.AutoFit()

# 6) Move the active selection as recorded in the sheet view.
$ws.Activate()
$ws.Range("G11").Select()
